# Update crypto price/volume data in the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on cells whose new values look numeric,
# so they remain text (matching the source data) instead of being
# parsed into floating point numbers by Excel.
$textRefs = @("D5", "D6", "D9", "D11", "D12", "D16", "D21", "D23", "D24", "D25", "D27", "D28", "D33", "D35", "D38", "D40", "D41", "D44", "D45", "D46", "D49")
foreach ($ref in $textRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "63.252.44"
$ws.Range("E2").Value = "  +5.11%  "
$ws.Range("D3").Value = "3.485.45"
$ws.Range("E3").Value = "  +4.45%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "409.22"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("D6").Value = "133.85"
$ws.Range("E6").Value = "  +20.40%  "
$ws.Range("D7").Value = "3.488.84"
$ws.Range("E7").Value = "  +4.92%  "
$ws.Range("E8").Value = "  +3.58%  "
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("E10").Value = "  +10.13%  "
$ws.Range("D11").Value = "0.131"
$ws.Range("E11").Value = "  +33.77%  "
$ws.Range("D12").Value = "43.62"
$ws.Range("E12").Value = "  +10.58%  "
$ws.Range("E13").Value = "  -0.57%  "
$ws.Range("D14").Value = "4.029.26"
$ws.Range("E14").Value = "  +4.02%  "
$ws.Range("E15").Value = "  +4.83%  "
$ws.Range("D16").Value = "20.33"
$ws.Range("E16").Value = "  +3.04%  "
$ws.Range("D17").Value = "3.494.24"
$ws.Range("E17").Value = "  +3.81%  "
$ws.Range("D18").Value = "63.126.26"
$ws.Range("E18").Value = "  +5.18%  "
$ws.Range("E19").Value = "  +1.33%  "
$ws.Range("E20").Value = "  +2.32%  "
$ws.Range("D21").Value = "0.0000139"
$ws.Range("E21").Value = "  +27.41%  "
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").Value = "83.03"
$ws.Range("E23").Value = "  +10.46%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").Value = "13.27"
$ws.Range("E24").Value = "  +0.89%  "
$ws.Range("D25").Value = "314.22"
$ws.Range("E25").Value = "  +4.76%  "
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("D27").Value = "30.57"
$ws.Range("E27").Value = "  +7.10%  "
$ws.Range("D28").Value = "8.29"
$ws.Range("E28").Value = "  +3.19%  "
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("E32").Value = "  +4.40%  "
$ws.Range("D33").Value = "44.19"
$ws.Range("E33").Value = "  +12.01%  "
$ws.Range("E34").Value = "  +3.15%  "
$ws.Range("D35").Value = "11.85"
$ws.Range("E35").Value = "  +3.57%  "
$ws.Range("E37").Value = "  -1.78%  "
$ws.Range("D38").Value = "52.64"
$ws.Range("E38").Value = "  +0.82%  "
$ws.Range("E39").Value = "  +6.59%  "
$ws.Range("D40").Value = "0.997"
$ws.Range("E40").Value = "  -0.28%  "
$ws.Range("D41").Value = "3.05"
$ws.Range("E41").Value = "  -0.98%  "
$ws.Range("E42").Value = "  +3.00%  "
$ws.Range("E43").Value = "  +4.36%  "
$ws.Range("D44").Value = "137.09"
$ws.Range("E44").Value = "  -0.45%  "
$ws.Range("D45").Value = "17.49"
$ws.Range("E45").Value = "  +3.82%  "
$ws.Range("D46").Value = "4.02"
$ws.Range("E46").Value = "  +2.51%  "
$ws.Range("E47").Value = "  -1.20%  "
$ws.Range("E48").Value = "  -0.90%  "
$ws.Range("D49").Value = "22.37"
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("D50").Value = "2.203.61"
$ws.Range("E50").Value = "  +0.54%  "
$ws.Range("D51").Value = "3.828.26"
$ws.Range("E51").Value = "  +3.92%  "

# Restore default (General) number format and Normal style so the
# cells keep looking like the rest of the untouched text cells.
foreach ($ref in $textRefs) {
    $ws.Range($ref).NumberFormat = "General"
    $ws.Range($ref).Style = "Normal"
}
